$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room for the new rows.
#    Current data occupies rows 2 (Morning), 3 (Afternoon), 4 (Night).
#    Target layout:
#      row1 : new header row (A..J)
#      row2 : Morning (unchanged data, restyled)
#      row3 : new split-off row (Selim / 0164576654)
#      row4 : new split-off row (Mosaddek / 0164576654)
#      row5 : Afternoon (old row 3 data)
#      row6 : Night (old row 4 data)
#    Row 1 is already blank, so we only need to insert two blank rows
#    right above the current row 3 - that pushes old row3->row5 and
#    old row4->row6 while leaving row1/row2 alone.
# ---------------------------------------------------------------------------
$ws.Rows("3:4").Insert()

# ---------------------------------------------------------------------------
# 2) New header row (row 1) - generic column headers A..J
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "A"
$ws.Range("B1").Value = "B"
$ws.Range("C1").Value = "C"
$ws.Range("D1").Value = "D"
$ws.Range("E1").Value = "E"
$ws.Range("F1").Value = "F"
$ws.Range("G1").Value = "G"
$ws.Range("H1").Value = "H"
$ws.Range("I1").Value = "I"
$ws.Range("J1").Value = "J"

# ---------------------------------------------------------------------------
# 3) New split-off rows 3 & 4 (only columns C & D are populated)
# ---------------------------------------------------------------------------
$ws.Range("C3").Value = "Selim"
$ws.Range("D3").Value = "0164576654"
$ws.Range("C4").Value = "Mosaddek"
$ws.Range("D4").Value = "0164576654"

# The row-insert carried column-level styles into columns F, H and J even
# though nothing is entered there for rows 3-4 in the target - drop that
# inherited formatting so those cells go back to being plain/empty.
$ws.Range("F3:F4").Style = "Normal"
$ws.Range("H3:H4").Style = "Normal"
$ws.Range("J3:J4").Style = "Normal"

# ---------------------------------------------------------------------------
# 4) Formatting - column widths
# ---------------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 18.7109375
$ws.Columns("B").ColumnWidth = 14.42578125
$ws.Columns("C").ColumnWidth = 11.5703125

# ---------------------------------------------------------------------------
# 5) Formatting - center the header cell A1 and the whole date/shift columns
# ---------------------------------------------------------------------------
$ws.Range("A1").VerticalAlignment = -4108
$ws.Range("A1").HorizontalAlignment = -4108

$ws.Range("B2:B4").VerticalAlignment = -4108
$ws.Range("B2:B4").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 6) Merge the date column across the whole shift block, and the shift-name
#    column across the Morning block (rows 2-4 share one "Morning" label).
# ---------------------------------------------------------------------------
$ws.Range("A2:A6").Merge()
$ws.Range("B2:B4").Merge()

# ---------------------------------------------------------------------------
# 7) Selection, matching the saved workbook's cursor position.
# ---------------------------------------------------------------------------
$ws.Range("D11").Select()
